$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.735.31"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "2.654.49"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'599.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").Value = "'168.57"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").Value = "2.654.58"
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "'0.367"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "3.137.86"
$ws.Range("E15").Value = "  -3.05%  "
$ws.Range("D16").Value = "'0.0000185"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "67.697.23"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "2.646.28"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").Value = "'12.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'7.91"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").Value = "'363.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("D23").Value = "'4.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.34%  "
$ws.Range("D24").Value = "'11.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.94%  "
$ws.Range("E25").Value = "  -4.49%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'70.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.08%  "
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'560.32"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.12%  "
$ws.Range("D32").Value = "'8.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("E33").Value = "  -4.28%  "
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'158.24"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'19.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").Value = "'0.374"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  -5.74%  "
$ws.Range("D46").Value = "'40.33"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("D47").Value = "0.0₆0300"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").Value = "'154.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").Value = "'3.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("E51").Value = "  -3.35%  "
